{"js": "// Rewrite the \"KEY ACHIEVEMENTS AND IMPACT\" > \"Impact\" bullet list from six\n// job-duty style bullets down to four impact-focused accomplishment bullets.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"KEY ACHIEVEMENTS AND IMPACT\" section heading, then the bullet\n// paragraphs that immediately follow its \"Impact\" sub-heading. Searching by\n// position (rather than a global text search/replace) is required because\n// two of the six legacy bullets are duplicated verbatim earlier in the\n// document (under PROFESSIONAL EXPERIENCE).\nlet sectionIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === \"KEY ACHIEVEMENTS AND IMPACT\") {\n    sectionIdx = i;\n    break;\n  }\n}\nif (sectionIdx === -1) {\n  throw new Error(\"Could not find 'KEY ACHIEVEMENTS AND IMPACT' heading\");\n}\n\n// The bullets begin after the heading + its \"Impact\" sub-heading paragraph.\nlet bulletStart = sectionIdx + 1;\nif (items[bulletStart] && items[bulletStart].text.trim() === \"Impact\") {\n  bulletStart += 1;\n}\n\n// Collect the run of bullet (\"\u2022\") paragraphs following the sub-heading.\nconst bulletIdxs = [];\nfor (let i = bulletStart; i < items.length && items[i].text.trim().startsWith(\"\u2022\"); i++) {\n  bulletIdxs.push(i);\n}\n\nconst oldBullets = [\n  \"\u2022 Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over $2 trillion\",\n  \"\u2022 Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy\",\n  \"\u2022 Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets\",\n  \"\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis\",\n  \"\u2022 Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\",\n  \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\",\n];\n\nconst newBullets = [\n  \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\",\n  \"\u2022 $4.7M savings enabled nonprofit access\",\n  \"\u2022 Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\",\n  \"\u2022 178% accuracy improvement in racial classification algorithms\",\n];\n\nif (bulletIdxs.length !== oldBullets.length) {\n  throw new Error(\n    \"Expected \" + oldBullets.length + \" bullets under KEY ACHIEVEMENTS, found \" + bulletIdxs.length\n  );\n}\nfor (let i = 0; i < bulletIdxs.length; i++) {\n  if (items[bulletIdxs[i]].text !== oldBullets[i]) {\n    throw new Error(\"Bullet \" + i + \" text mismatch: \" + items[bulletIdxs[i]].text);\n  }\n}\n\n// Replace the text of the first four bullet paragraphs with the new\n// accomplishment statements, then delete the trailing two paragraphs that no\n// longer have a replacement (6 bullets -> 4 bullets).\nfor (let i = 0; i < newBullets.length; i++) {\n  items[bulletIdxs[i]].insertText(newBullets[i], Word.InsertLocation.replace);\n}\nfor (let i = newBullets.length; i < bulletIdxs.length; i++) {\n  items[bulletIdxs[i]].delete();\n}\n\nawait context.sync();\n", "ps1": "# Rewrite the \"KEY ACHIEVEMENTS AND IMPACT\" > \"Impact\" bullet list from six\n# job-duty style bullets down to four impact-focused accomplishment bullets.\n$d = $word.ActiveDocument\n\n# Locate the \"KEY ACHIEVEMENTS AND IMPACT\" section heading, then the bullet\n# paragraphs that immediately follow its \"Impact\" sub-heading. Walking\n# paragraphs by position (rather than a global Find/Replace) is required\n# because two of the six legacy bullets are duplicated verbatim earlier in\n# the document (under PROFESSIONAL EXPERIENCE).\n$sectionIdx = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text.Trim()\n    if ($t -eq \"KEY ACHIEVEMENTS AND IMPACT\") {\n        $sectionIdx = $i\n        break\n    }\n}\nif ($sectionIdx -eq -1) {\n    throw \"Could not find 'KEY ACHIEVEMENTS AND IMPACT' heading\"\n}\n\n# The bullets begin after the heading + its \"Impact\" sub-heading paragraph.\n$bulletStart = $sectionIdx + 1\nif ($d.Paragraphs($bulletStart).Range.Text.Trim() -eq \"Impact\") {\n    $bulletStart = $bulletStart + 1\n}\n\n# Collect the run of bullet (\"\u2022\") paragraphs following the sub-heading.\n$bulletIdxs = @()\n$i = $bulletStart\nwhile ($i -le $d.Paragraphs.Count) {\n    $t = $d.Paragraphs($i).Range.Text.Trim()\n    if ($t.StartsWith(\"\u2022\")) {\n        $bulletIdxs += $i\n        $i = $i + 1\n    } else {\n        break\n    }\n}\n\n$oldBullets = @(\n    \"\u2022 Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over `$2 trillion\",\n    \"\u2022 Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy\",\n    \"\u2022 Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets\",\n    \"\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis\",\n    \"\u2022 Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\",\n    \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\"\n)\n\n$newBullets = @(\n    \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\",\n    \"\u2022 `$4.7M savings enabled nonprofit access\",\n    \"\u2022 Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\",\n    \"\u2022 178% accuracy improvement in racial classification algorithms\"\n)\n\nif ($bulletIdxs.Count -ne $oldBullets.Count) {\n    throw \"Expected $($oldBullets.Count) bullets under KEY ACHIEVEMENTS, found $($bulletIdxs.Count)\"\n}\nfor ($k = 0; $k -lt $bulletIdxs.Count; $k++) {\n    $actual = $d.Paragraphs($bulletIdxs[$k]).Range.Text.Trim()\n    if ($actual -ne $oldBullets[$k]) {\n        throw \"Bullet $k text mismatch: $actual\"\n    }\n}\n\n# Replace the text of the first four bullet paragraphs with the new\n# accomplishment statements. Assigning `.Range.Text` (with no trailing `\\r`)\n# replaces just the paragraph's text run and leaves its own paragraph mark\n# (and paragraph count) untouched; appending `\\r` would insert an extra\n# paragraph break instead of a plain replace.\nfor ($k = 0; $k -lt $newBullets.Count; $k++) {\n    $idx = $bulletIdxs[$k]\n    $d.Paragraphs($idx).Range.Text = $newBullets[$k]\n}\n\n# Delete the trailing two paragraphs that no longer have a replacement\n# (6 bullets -> 4 bullets). Walk backwards so earlier indices stay valid.\nfor ($k = $bulletIdxs.Count - 1; $k -ge $newBullets.Count; $k--) {\n    $idx = $bulletIdxs[$k]\n    $d.Paragraphs($idx).Range.Delete()\n}\n"}
